{"js": "// Insert a new bold line at the very top of the document with the\n// \"uvicorn\" run command, followed by two blank (bold) lines, right\n// before the existing \"AI Prompt Book v3\" title paragraph.\n//\n// We use insertOoxml() so the inserted markup matches exactly how Word\n// itself would emit it after a spell-check pass (proofErr spellStart /\n// spellEnd markers bracketing the two \"misspelled\" tokens \"uvicorn\" and\n// \"app:app\", each sub-run carrying its own bold run properties).\n\nconst body = context.document.body;\nconst startRange = body.getRange(\"Start\");\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>uvicorn</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>app:app</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> --host 0.0.0.0 --port $PORT</w:t></w:r>' +\n            '</w:p>' +\n            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>' +\n            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nstartRange.insertOoxml(ooxml, \"Before\");\nawait context.sync();\n", "ps1": "# Insert a new bold line at the very top of the document with the\n# \"uvicorn\" run command, followed by two blank (bold) lines, right\n# before the existing \"AI Prompt Book v3\" title paragraph.\n#\n# We build the exact OOXML Word itself would emit after a spell-check\n# pass (proofErr spellStart / spellEnd markers bracketing the two\n# \"misspelled\" tokens \"uvicorn\" and \"app:app\", each sub-run carrying its\n# own bold run properties) and inject it with Range.InsertXML so the\n# resulting markup matches precisely.\n\n$d = $word.ActiveDocument\n\n$target = $d.Paragraphs.Item(1).Range\n$target.Collapse(1)  # wdCollapseStart\n\n$ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>uvicorn</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>app:app</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\"> --host 0.0.0.0 --port $PORT</w:t></w:r>' +\n            '</w:p>' +\n            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>' +\n            '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n[void]$target.InsertXML($ooxml)\n"}
